# Generate Report for Handback
# Updates the handback-status workbook so that the two e2e file rows that
# used to reference:
#   6996298e-972a-422f-874d-04f8a2062c52.md
#   dd5a7bb6-83dc-49c5-b38f-53e74a0b6ef0.md
# now reference:
#   3914d24b-8b11-4bc1-8ddc-84a65dd1ee83.md
#   ffffb474439f-8e90-4131-8340-90fa51f53248.md
# and refreshes the correspond handoff/handback xlf names + timestamps that
# go with the new handback run, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldFile1 = "6996298e-972a-422f-874d-04f8a2062c52.md"
$newFile1 = "3914d24b-8b11-4bc1-8ddc-84a65dd1ee83.md"
$oldFile2 = "dd5a7bb6-83dc-49c5-b38f-53e74a0b6ef0.md"
$newFile2 = "ffffb474439f-8e90-4131-8340-90fa51f53248.md"

$oldXlfZh1 = "6996298e-972a-422f-874d-04f8a2062c52.a98dc62e5417b615d9831d9716010847cc3c2164.zh-cn.xlf"
$oldXlfZh2 = "dd5a7bb6-83dc-49c5-b38f-53e74a0b6ef0.022f0846ef4cd392eb68e416f92ebbffcf23b4a4.zh-cn.xlf"
$newXlfZh  = "3914d24b-8b11-4bc1-8ddc-84a65dd1ee83.d97f4267a936c10b3e28a3a56e067270e7c460a3.zh-cn.xlf"

$oldXlfDe1 = "6996298e-972a-422f-874d-04f8a2062c52.a98dc62e5417b615d9831d9716010847cc3c2164.de-de.xlf"
$oldXlfDe2 = "dd5a7bb6-83dc-49c5-b38f-53e74a0b6ef0.022f0846ef4cd392eb68e416f92ebbffcf23b4a4.de-de.xlf"
$newXlfDe  = "3914d24b-8b11-4bc1-8ddc-84a65dd1ee83.d97f4267a936c10b3e28a3a56e067270e7c460a3.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = "e2e\" + $newFile1
$wsOverview.Range("G2").Value = "2016-08-18 19:07:18"

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = "e2e\" + $newFile2
$wsOverview.Range("G3").Value = "2016-08-18 19:07:18"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce9f262e75a5942f24066cff7a650248e188a9e0/e2e/" + $oldFile1, "", "", "e2e\" + $newFile1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce9f262e75a5942f24066cff7a650248e188a9e0/e2e/" + $oldFile2, "", "", "e2e\" + $newFile2)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("G2").Value = $newXlfZh
$wsZhCn.Range("H2").Value = "2016-08-18 19:07:12"
$wsZhCn.Range("I2").Value = $newFile1
$wsZhCn.Range("J2").Value = $newXlfZh
$wsZhCn.Range("K2").Value = "2016-08-18 19:07:41"

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("G3").Value = $newXlfZh
$wsZhCn.Range("H3").Value = "2016-08-18 19:07:12"
$wsZhCn.Range("I3").Value = $newFile2
$wsZhCn.Range("J3").Value = $newXlfZh
$wsZhCn.Range("K3").Value = "2016-08-18 19:07:41"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce9f262e75a5942f24066cff7a650248e188a9e0/e2e/" + $oldFile1, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/abcb9a6cfaae25d79d8cb65215ad8c575b9aa754/e2e/" + $oldFile1, "", "", $newFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce9f262e75a5942f24066cff7a650248e188a9e0/e2e/" + $oldFile2, "", "", $newFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/abcb9a6cfaae25d79d8cb65215ad8c575b9aa754/e2e/" + $oldFile2, "", "", $newFile2)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("G2").Value = $newXlfDe
$wsDeDe.Range("H2").Value = "2016-08-18 19:07:18"
$wsDeDe.Range("I2").Value = $newFile1
$wsDeDe.Range("J2").Value = $newXlfDe
$wsDeDe.Range("K2").Value = "2016-08-18 19:07:49"

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("G3").Value = $newXlfDe
$wsDeDe.Range("H3").Value = "2016-08-18 19:07:18"
$wsDeDe.Range("I3").Value = $newFile2
$wsDeDe.Range("J3").Value = $newXlfDe
$wsDeDe.Range("K3").Value = "2016-08-18 19:07:49"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce9f262e75a5942f24066cff7a650248e188a9e0/e2e/" + $oldFile1, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4767460c85cd2bdcf599033861fc233b0d74bbda/e2e/" + $oldFile1, "", "", $newFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce9f262e75a5942f24066cff7a650248e188a9e0/e2e/" + $oldFile2, "", "", $newFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4767460c85cd2bdcf599033861fc233b0d74bbda/e2e/" + $oldFile2, "", "", $newFile2)

Write-Output "Done updating handback status workbook."
